# Updated symbol list on Thu Dec 29 09:54:09 UTC 2022 with GitHub Actions
# Refreshes the cryptocurrency Price / Volume(1h) rank columns on Sheet1.
# Note: Price values (column D) look numeric but must stay stored as text
# (matching the original inlineStr cells), so they are entered with a
# leading apostrophe to force Excel to keep them as text instead of
# auto-converting to a Double (which would lose exact formatting such as
# trailing/leading zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.78"
$ws.Range("D5").Value = "'0.05738"
$ws.Range("D6").Value = "'6.488"
$ws.Range("D7").Value = "'3.165"
$ws.Range("D8").Value = "'0.8149"
$ws.Range("D9").Value = "'0.8569"
$ws.Range("D10").Value = "'0.1376"
$ws.Range("D12").Value = "'0.03178"
$ws.Range("D13").Value = "'0.02875"
$ws.Range("D14").Value = "'0.09344"
$ws.Range("D15").Value = "'3.817"
$ws.Range("D16").Value = "'0.001525"
$ws.Range("D17").Value = "'0.04699"
$ws.Range("D18").Value = "'0.0005996"
$ws.Range("D19").Value = "'0.006229"
$ws.Range("D20").Value = "'0.001235"
$ws.Range("D21").Value = "'0.004792"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00008493"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.532"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.153"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3201"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1338"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("B27").Value = "ZBToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D27").Value = "'0.1327"
$ws.Range("E27").Value = "26ZBTokenZB"
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D28").Value = "'0.0002329"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("B29").Value = "Spectre.aiUtilityToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("E29").Value = "28Spectre.aiUtilityTokenSXUT"
$ws.Range("B30").Value = "LegolasExchange"
$ws.Range("C30").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("E30").Value = "29LegolasExchangeLGO"
$ws.Range("B31").Value = "BitZToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("E31").Value = "30BitZTokenBZ"
$ws.Range("B32").Value = "Birake"
$ws.Range("C32").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("E32").Value = "31BirakeBIR"
$ws.Range("D40").Value = "'0.03694"
$ws.Range("D41").Value = "'0.006394"
$ws.Range("D42").Value = "'0.1055"
$ws.Range("D43").Value = "'0.002252"
$ws.Range("D44").Value = "'0.007798"
$ws.Range("D45").Value = "'0.00005493"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.3997"
$ws.Range("D48").Value = "'0.002508"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.0001998"
